# One-click update from Van Paper 07:26 AM on 2025-10-23
# Inserts a new leaderboard row for customer "WOODBURY ICE" (0008347),
# salesperson "Norman, Ryan M" (040), pushing the former row 22
# ("HOLY FAMILY MARONITE CHURCH") down to row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 22 - this shifts row 22
# ("HOLY FAMILY MARONITE CHURCH" ...) down to row 23 and inherits the
# formatting of the row above it (row 21), which matches the target
# style pattern for the new row.
$ws.Rows.Item(22).Insert()

# Match the row height used by every other data row.
$ws.Rows.Item(22).RowHeight = 13.05

# Populate the newly-inserted row 22 with the new leaderboard entry.
$ws.Cells.Item(22, 1).Value = "WOODBURY ICE"
$ws.Cells.Item(22, 2).Value = "Norman, Ryan M"
$ws.Cells.Item(22, 3).Value = "040"
$ws.Cells.Item(22, 5).Value = "0008347"

# Column F stays blank for this row (as it does for every other row) -
# copy the untouched blank cell from the row above so the cell is
# materialized the same way the rest of the sheet's blank F cells are.
$ws.Cells.Item(21, 6).Copy($ws.Cells.Item(22, 6))
